$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: permute full row data (columns B:AC) among swapped/rotated match rows ---

# cycle: [13, 14]
$row13 = $ws.Range("B13:AC13").Value2
$row14 = $ws.Range("B14:AC14").Value2
$ws.Range("B13:AC13").Value = $row14
$ws.Range("B14:AC14").Value = $row13

# cycle: [53, 54]
$row53 = $ws.Range("B53:AC53").Value2
$row54 = $ws.Range("B54:AC54").Value2
$ws.Range("B53:AC53").Value = $row54
$ws.Range("B54:AC54").Value = $row53

# cycle: [87, 89]
$row87 = $ws.Range("B87:AC87").Value2
$row89 = $ws.Range("B89:AC89").Value2
$ws.Range("B87:AC87").Value = $row89
$ws.Range("B89:AC89").Value = $row87

# cycle: [131, 133]
$row131 = $ws.Range("B131:AC131").Value2
$row133 = $ws.Range("B133:AC133").Value2
$ws.Range("B131:AC131").Value = $row133
$ws.Range("B133:AC133").Value = $row131

# cycle: [144, 146, 147]
$row144 = $ws.Range("B144:AC144").Value2
$row146 = $ws.Range("B146:AC146").Value2
$row147 = $ws.Range("B147:AC147").Value2
$ws.Range("B144:AC144").Value = $row146
$ws.Range("B146:AC146").Value = $row147
$ws.Range("B147:AC147").Value = $row144

# cycle: [158, 159]
$row158 = $ws.Range("B158:AC158").Value2
$row159 = $ws.Range("B159:AC159").Value2
$ws.Range("B158:AC158").Value = $row159
$ws.Range("B159:AC159").Value = $row158

# cycle: [172, 174, 175, 176, 173]
$row172 = $ws.Range("B172:AC172").Value2
$row174 = $ws.Range("B174:AC174").Value2
$row175 = $ws.Range("B175:AC175").Value2
$row176 = $ws.Range("B176:AC176").Value2
$row173 = $ws.Range("B173:AC173").Value2
$ws.Range("B172:AC172").Value = $row174
$ws.Range("B174:AC174").Value = $row175
$ws.Range("B175:AC175").Value = $row176
$ws.Range("B176:AC176").Value = $row173
$ws.Range("B173:AC173").Value = $row172

# cycle: [186, 187]
$row186 = $ws.Range("B186:AC186").Value2
$row187 = $ws.Range("B187:AC187").Value2
$ws.Range("B186:AC186").Value = $row187
$ws.Range("B187:AC187").Value = $row186

# cycle: [189, 190, 191]
$row189 = $ws.Range("B189:AC189").Value2
$row190 = $ws.Range("B190:AC190").Value2
$row191 = $ws.Range("B191:AC191").Value2
$ws.Range("B189:AC189").Value = $row190
$ws.Range("B190:AC190").Value = $row191
$ws.Range("B191:AC191").Value = $row189

# cycle: [200, 201]
$row200 = $ws.Range("B200:AC200").Value2
$row201 = $ws.Range("B201:AC201").Value2
$ws.Range("B200:AC200").Value = $row201
$ws.Range("B201:AC201").Value = $row200

# cycle: [202, 203]
$row202 = $ws.Range("B202:AC202").Value2
$row203 = $ws.Range("B203:AC203").Value2
$ws.Range("B202:AC202").Value = $row203
$ws.Range("B203:AC203").Value = $row202

# --- Step 2: flip FK Indija <-> FK Tekstilac Odzaci team-name text (shared-string swap effect) ---
$ws.Range("G2").Value = "FK Tekstilac Odzaci"
$ws.Range("G9").Value = "FK Indija"
$ws.Range("G20").Value = "FK Indija"
$ws.Range("G23").Value = "FK Tekstilac Odzaci"
$ws.Range("F28").Value = "FK Tekstilac Odzaci"
$ws.Range("F29").Value = "FK Indija"
$ws.Range("G35").Value = "FK Tekstilac Odzaci"
$ws.Range("G39").Value = "FK Indija"
$ws.Range("F43").Value = "FK Tekstilac Odzaci"
$ws.Range("F46").Value = "FK Indija"
$ws.Range("G51").Value = "FK Tekstilac Odzaci"
$ws.Range("F58").Value = "FK Tekstilac Odzaci"
$ws.Range("F63").Value = "FK Indija"
$ws.Range("F70").Value = "FK Tekstilac Odzaci"
$ws.Range("G70").Value = "FK Indija"
$ws.Range("G75").Value = "FK Tekstilac Odzaci"
$ws.Range("F77").Value = "FK Indija"
$ws.Range("F82").Value = "FK Tekstilac Odzaci"
$ws.Range("G90").Value = "FK Tekstilac Odzaci"
$ws.Range("F97").Value = "FK Indija"
$ws.Range("G102").Value = "FK Indija"
$ws.Range("F104").Value = "FK Tekstilac Odzaci"
$ws.Range("F107").Value = "FK Indija"
$ws.Range("G112").Value = "FK Tekstilac Odzaci"
$ws.Range("F114").Value = "FK Tekstilac Odzaci"
$ws.Range("F116").Value = "FK Indija"
$ws.Range("G124").Value = "FK Indija"
$ws.Range("G132").Value = "FK Indija"
$ws.Range("F138").Value = "FK Tekstilac Odzaci"
$ws.Range("F139").Value = "FK Indija"
$ws.Range("G143").Value = "FK Tekstilac Odzaci"
$ws.Range("G145").Value = "FK Indija"
$ws.Range("F153").Value = "FK Tekstilac Odzaci"
$ws.Range("G161").Value = "FK Tekstilac Odzaci"
$ws.Range("G165").Value = "FK Indija"
$ws.Range("G168").Value = "FK Tekstilac Odzaci"
$ws.Range("F170").Value = "FK Indija"
$ws.Range("G179").Value = "FK Tekstilac Odzaci"
$ws.Range("G180").Value = "FK Indija"
$ws.Range("G205").Value = "FK Indija"

# --- Step 3: delete trailing row 208 (future fixture removed) ---
$ws.Rows.Item(208).Delete()

Write-Host "edit applied"
